$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("B3").Value = 79245
$ws.Range("B4").Value = 79245
$ws.Range("B5").Value = 79245
$ws.Range("B6").Value = 79864
